$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 351; this pushes the former rows
# 351..432 down to 352..433 (all their values/styles shift down with them).
$ws.Rows.Item(351).Insert()

# Populate the newly inserted (now-blank) row 351 with the new record.
$ws.Range("A351").Value = 10
$ws.Range("B351").Value = "Vega Modelo de Temuco"
$ws.Range("C351").Value = "La Araucanía"
$ws.Range("D351").Value = 44511
$ws.Range("E351").Value = 9
$ws.Range("F351").Value = "Fruta"
$ws.Range("G351").Value = 100102
$ws.Range("H351").Value = "Cítricos"
$ws.Range("I351").Value = 100102004
$ws.Range("J351").Value = "Mandarina"
$ws.Range("K351").Value = "Clementina"
$ws.Range("L351").Value = "Primera"
$ws.Range("M351").Value = 115
$ws.Range("N351").Value = 9000
$ws.Range("O351").Value = 10000
$ws.Range("P351").Value = 9565
$ws.Range("Q351").Value = "$/bandeja 18 kilos"
$ws.Range("R351").Value = "Región de O'Higgins"
$ws.Range("S351").Value = 531
$ws.Range("T351").Value = 18
